$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.441829000000001
$ws.Range("H2").Value = 16.325487
$ws.Range("I2").Value = 0.5729403216841985
$ws.Range("J2").Value = 0.5729403216841985
$ws.Range("M2").Value = 57.65261933333333
$ws.Range("N2").Value = 172.957858
$ws.Range("O2").Value = 0.6817060950001529
$ws.Range("P2").Value = 0.6817060950001529
$ws.Range("Q2").Value = 313.735695814094
$ws.Range("R2").Value = 2823.621262326846
$ws.Range("S2").Value = 0.3905769093634664
$ws.Range("T2").Value = 0.3905769093634664
$ws.Range("G3").Value = 5.441829000000001
$ws.Range("H3").Value = 16.325487
$ws.Range("I3").Value = 0.5729403216841985
$ws.Range("J3").Value = 0.5729403216841985
$ws.Range("O3").Value = 0.1019529789289588
$ws.Range("P3").Value = 0.1019529789289588
$ws.Range("Q3").Value = 46.92093707126001
$ws.Range("R3").Value = 422.28843364134
$ws.Range("S3").Value = 0.05841297254421998
$ws.Range("T3").Value = 0.05841297254421998
$ws.Range("G4").Value = 5.441829000000001
$ws.Range("H4").Value = 16.325487
$ws.Range("I4").Value = 0.5729403216841985
$ws.Range("J4").Value = 0.5729403216841985
$ws.Range("M4").Value = 2.790736
$ws.Range("N4").Value = 8.372208000000001
$ws.Range("O4").Value = 0.0329987043561157
$ws.Range("P4").Value = 0.0329987043561157
$ws.Range("Q4").Value = 15.186708096144
$ws.Range("R4").Value = 136.680372865296
$ws.Range("S4").Value = 0.01890628828895469
$ws.Range("T4").Value = 0.01890628828895469
$ws.Range("G5").Value = 5.441829000000001
$ws.Range("H5").Value = 16.325487
$ws.Range("I5").Value = 0.5729403216841985
$ws.Range("J5").Value = 0.5729403216841985
$ws.Range("M5").Value = 15.50544933333333
$ws.Range("N5").Value = 46.516348
$ws.Range("O5").Value = 0.1833422217147727
$ws.Range("P5").Value = 0.1833422217147727
$ws.Range("Q5").Value = 84.37800384016401
$ws.Range("R5").Value = 759.4020345614762
$ws.Range("S5").Value = 0.1050441514875575
$ws.Range("T5").Value = 0.1050441514875575
$ws.Range("I6").Value = 0.2716201486343598
$ws.Range("J6").Value = 0.2716201486343598
$ws.Range("M6").Value = 57.65261933333333
$ws.Range("N6").Value = 172.957858
$ws.Range("O6").Value = 0.6817060950001529
$ws.Range("P6").Value = 0.6817060950001529
$ws.Range("Q6").Value = 148.736147734248
$ws.Range("R6").Value = 1338.625329608232
$ws.Range("S6").Value = 0.1851651108488905
$ws.Range("T6").Value = 0.1851651108488905
$ws.Range("I7").Value = 0.2716201486343598
$ws.Range("J7").Value = 0.2716201486343598
$ws.Range("O7").Value = 0.1019529789289588
$ws.Range("P7").Value = 0.1019529789289588
$ws.Range("S7").Value = 0.02769248329039955
$ws.Range("T7").Value = 0.02769248329039955
$ws.Range("I8").Value = 0.2716201486343598
$ws.Range("J8").Value = 0.2716201486343598
$ws.Range("M8").Value = 2.790736
$ws.Range("N8").Value = 8.372208000000001
$ws.Range("O8").Value = 0.0329987043561157
$ws.Range("P8").Value = 0.0329987043561157
$ws.Range("Q8").Value = 7.199730502848
$ws.Range("R8").Value = 64.797574525632
$ws.Range("S8").Value = 0.008963112981949442
$ws.Range("T8").Value = 0.008963112981949441
$ws.Range("I9").Value = 0.2716201486343598
$ws.Range("J9").Value = 0.2716201486343598
$ws.Range("M9").Value = 15.50544933333333
$ws.Range("N9").Value = 46.516348
$ws.Range("O9").Value = 0.1833422217147727
$ws.Range("P9").Value = 0.1833422217147727
$ws.Range("Q9").Value = 40.00201256068799
$ws.Range("R9").Value = 360.018113046192
$ws.Range("S9").Value = 0.0497994415131203
$ws.Range("T9").Value = 0.0497994415131203
$ws.Range("G10").Value = 1.476376
$ws.Range("H10").Value = 4.429128
$ws.Range("I10").Value = 0.1554395296814417
$ws.Range("J10").Value = 0.1554395296814417
$ws.Range("M10").Value = 57.65261933333333
$ws.Range("N10").Value = 172.957858
$ws.Range("O10").Value = 0.6817060950001529
$ws.Range("P10").Value = 0.6817060950001529
$ws.Range("Q10").Value = 85.11694352086933
$ws.Range("R10").Value = 766.052491687824
$ws.Range("S10").Value = 0.105964074787796
$ws.Range("T10").Value = 0.105964074787796
$ws.Range("G11").Value = 1.476376
$ws.Range("H11").Value = 4.429128
$ws.Range("I11").Value = 0.1554395296814417
$ws.Range("J11").Value = 0.1554395296814417
$ws.Range("O11").Value = 0.1019529789289588
$ws.Range("P11").Value = 0.1019529789289588
$ws.Range("Q11").Value = 12.72971741477333
$ws.Range("R11").Value = 114.56745673296
$ws.Range("S11").Value = 0.01584752309433929
$ws.Range("T11").Value = 0.0158475230943393
$ws.Range("G12").Value = 1.476376
$ws.Range("H12").Value = 4.429128
$ws.Range("I12").Value = 0.1554395296814417
$ws.Range("J12").Value = 0.1554395296814417
$ws.Range("M12").Value = 2.790736
$ws.Range("N12").Value = 8.372208000000001
$ws.Range("O12").Value = 0.0329987043561157
$ws.Range("P12").Value = 0.0329987043561157
$ws.Range("Q12").Value = 4.120175652736001
$ws.Range("R12").Value = 37.081580874624
$ws.Range("S12").Value = 0.005129303085211565
$ws.Range("T12").Value = 0.005129303085211565
$ws.Range("G13").Value = 1.476376
$ws.Range("H13").Value = 4.429128
$ws.Range("I13").Value = 0.1554395296814417
$ws.Range("J13").Value = 0.1554395296814417
$ws.Range("M13").Value = 15.50544933333333
$ws.Range("N13").Value = 46.516348
$ws.Range("O13").Value = 0.1833422217147727
$ws.Range("P13").Value = 0.1833422217147727
$ws.Range("Q13").Value = 22.89187326494934
$ws.Range("R13").Value = 206.026859384544
$ws.Range("S13").Value = 0.02849862871409487
$ws.Range("T13").Value = 0.02849862871409488
